$d = $word.ActiveDocument

$replacements = @(
    @("2023-09-21 Thursday", "2023-09-22 Friday"),
    @("27×97=2619", "60×62=3720"),
    @("32×11=352", "22×30=660"),
    @("77×92=7084", "67×74=4958"),
    @("85×86=7310", "50×58=2900"),
    @("31×66=2046", "92×52=4784"),
    @("14×60=840", "73×98=7154"),
    @("79×64=5056", "55×73=4015"),
    @("82×56=4592", "37×93=3441"),
    @("99×51=5049", "61×68=4148"),
    @("75×11=825", "82×54=4428"),
    @("36×59=2124", "17×31=527"),
    @("13×48=624", "59×51=3009"),
    @("47×71=3337", "84×55=4620"),
    @("18×64=1152", "45×94=4230"),
    @("70×98=6860", "35×83=2905"),
    @("43×41=1763", "82×37=3034"),
    @("63×45=2835", "32×78=2496"),
    @("13×17=221", "60×98=5880"),
    @("12×17=204", "14×67=938"),
    @("30×76=2280", "98×43=4214"),
    @("21×63=1323", "55×97=5335"),
    @("70×52=3640", "72×43=3096"),
    @("64×81=5184", "87×13=1131"),
    @("53×26=1378", "89×70=6230"),
    @("74×76=5624", "71×93=6603")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
